# Market-data refresh for the Seraph profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# on the affected Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# to match freshly pulled marketboard data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 6838  # was 5946.6665
$ws.Range("I43").Value = 3547.5  # was 3136
$ws.Range("K43").Value = 3547.5  # was 3136
$ws.Range("M43").Value = -3478.5  # was -3067

# Row 62
$ws.Range("H62").Value = 0  # was 3000
$ws.Range("I62").Value = 0  # was 3000
$ws.Range("K62").Value = 0  # was 3000
$ws.Range("M62").ClearContents()  # was -2376

# Row 65
$ws.Range("H65").Value = 0  # was 3000
$ws.Range("I65").Value = 0  # was 3000
$ws.Range("K65").Value = 0  # was 15000
$ws.Range("M65").ClearContents()  # was -11880

# Row 80
$ws.Range("H80").Value = 407.18182  # was 377.33334
$ws.Range("J80").Value = 487  # was 432.25
$ws.Range("L80").Value = 1461  # was 1296.75
$ws.Range("N80").Value = -3457  # was -3292.75

# Row 83
$ws.Range("H83").Value = 407.18182  # was 377.33334
$ws.Range("J83").Value = 487  # was 432.25
$ws.Range("L83").Value = 4383  # was 3890.25
$ws.Range("N83").Value = -14367  # was -13874.25

# Row 98
$ws.Range("H98").Value = 2256.6667  # was 4262.6665
$ws.Range("I98").Value = 2256.6667  # was 3314
$ws.Range("J98").Value = 0  # was 9006
$ws.Range("K98").Value = 2256.6667  # was 3314
$ws.Range("L98").Value = 0  # was 9006
$ws.Range("M98").Value = -758.6667000000002  # was -1816
$ws.Range("N98").ClearContents()  # was -12002

# Row 118
$ws.Range("H118").Value = 537.8  # was 668
$ws.Range("I118").Value = 537.8  # was 585
$ws.Range("J118").Value = 0  # was 1000
$ws.Range("K118").Value = 1613.4  # was 1755
$ws.Range("L118").Value = 0  # was 3000
$ws.Range("M118").Value = 43.60000000000014  # was -98
$ws.Range("N118").ClearContents()  # was -6314

# Row 122
$ws.Range("H122").Value = 2256.6667  # was 4262.6665
$ws.Range("I122").Value = 2256.6667  # was 3314
$ws.Range("J122").Value = 0  # was 9006
$ws.Range("K122").Value = 6770.000100000001  # was 9942
$ws.Range("L122").Value = 0  # was 27018
$ws.Range("M122").Value = -4320.000100000001  # was -7492
$ws.Range("N122").ClearContents()  # was -31918

# Row 138
$ws.Range("H138").Value = 1893.8518  # was 1942.72
$ws.Range("I138").Value = 1778.72  # was 1815.3334
$ws.Range("J138").Value = 3333  # was 5000
$ws.Range("K138").Value = 5336.16  # was 5446.0002
$ws.Range("L138").Value = 9999  # was 15000
$ws.Range("M138").Value = -196.1599999999999  # was -306.0002000000004
$ws.Range("N138").Value = -20279  # was -25280

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3342  # was 3173
$ws.Range("I45").Value = 1997  # was 1994.6666
$ws.Range("K45").Value = 1997  # was 1994.6666
$ws.Range("M45").Value = -1620  # was -1617.6666

# Row 122
$ws.Range("H122").Value = 15465.75  # was 18100
$ws.Range("I122").Value = 8241  # was 9940
$ws.Range("K122").Value = 24723  # was 29820
$ws.Range("M122").Value = -22273  # was -27370

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 19283.166  # was 19598.834
$ws.Range("J35").Value = 19283.166  # was 19598.834
$ws.Range("L35").Value = 19283.166  # was 19598.834
$ws.Range("N35").Value = -19903.166  # was -20218.834

# Row 94
$ws.Range("H94").Value = 1055.0968  # was 1058.3549
$ws.Range("I94").Value = 547.7917  # was 552
$ws.Range("K94").Value = 547.7917  # was 552
$ws.Range("M94").Value = -96.79169999999999  # was -101

# Row 105
$ws.Range("H105").Value = 5212332  # was 4633295.5
$ws.Range("I105").Value = 20838828  # was 16671265
$ws.Range("J105").Value = 3499.75  # was 3307.4614
$ws.Range("K105").Value = 20838828  # was 16671265
$ws.Range("L105").Value = 3499.75  # was 3307.4614
$ws.Range("M105").Value = -20837081  # was -16669518
$ws.Range("N105").Value = -6993.75  # was -6801.4614

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 29415778  # was 125753750
$ws.Range("I16").Value = 35717588  # was 250002500
$ws.Range("J16").Value = 7333.3335  # was 1505000.5
$ws.Range("K16").Value = 35717588  # was 250002500
$ws.Range("L16").Value = 7333.3335  # was 1505000.5
$ws.Range("M16").Value = -35717301  # was -250002213
$ws.Range("N16").Value = -7907.3335  # was -1505574.5

# Row 113
$ws.Range("H113").Value = 29415778  # was 125753750
$ws.Range("I113").Value = 35717588  # was 250002500
$ws.Range("J113").Value = 7333.3335  # was 1505000.5
$ws.Range("K113").Value = 35717588  # was 250002500
$ws.Range("L113").Value = 7333.3335  # was 1505000.5
$ws.Range("M113").Value = -35715418  # was -250000330
$ws.Range("N113").Value = -11673.3335  # was -1509340.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 66911.336  # was 66908.2
$ws.Range("I2").Value = 166720.17  # was 142903.42
$ws.Range("J2").Value = 372.1111  # was 412.375
$ws.Range("K2").Value = 1000321.02  # was 857420.52
$ws.Range("L2").Value = 2232.6666  # was 2474.25
$ws.Range("M2").Value = -1000208.02  # was -857307.52
$ws.Range("N2").Value = -2458.6666  # was -2700.25

# Row 4
$ws.Range("H4").Value = 591.55554  # was 607.7059
$ws.Range("I4").Value = 524.2273  # was 559.4138
$ws.Range("K4").Value = 1572.6819  # was 1678.2414
$ws.Range("M4").Value = -1460.6819  # was -1566.2414

# Row 9
$ws.Range("H9").Value = 350  # was 7000175
$ws.Range("J9").Value = 0  # was 14000000
$ws.Range("L9").Value = 0  # was 42000000
$ws.Range("N9").ClearContents()  # was -42000448

# Row 97
$ws.Range("H97").Value = 4673.4  # was 4792
$ws.Range("J97").Value = 4726  # was 4989.5
$ws.Range("L97").Value = 14178  # was 14968.5
$ws.Range("N97").Value = -15170  # was -15960.5

# Row 98
$ws.Range("H98").Value = 3963.5557  # was 4012.4285
$ws.Range("J98").Value = 3813.1667  # was 3823.5
$ws.Range("L98").Value = 11439.5001  # was 11470.5
$ws.Range("N98").Value = -14435.5001  # was -14466.5

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 2754376  # was 2755001
$ws.Range("I18").Value = 3667501.8  # was 5500002.5
$ws.Range("J18").Value = 14999  # was 9999.5
$ws.Range("K18").Value = 3667501.8  # was 5500002.5
$ws.Range("L18").Value = 14999  # was 9999.5
$ws.Range("M18").Value = -3667208.8  # was -5499709.5
$ws.Range("N18").Value = -15585  # was -10585.5

# Row 29
$ws.Range("H29").Value = 0  # was 1500
$ws.Range("J29").Value = 0  # was 1500
$ws.Range("L29").Value = 0  # was 1500
$ws.Range("N29").ClearContents()  # was -2080

# Row 80
$ws.Range("H80").Value = 4147.7  # was 4322.65
$ws.Range("I80").Value = 3424.9  # was 3694.3333
$ws.Range("J80").Value = 4870.5  # was 4836.727
$ws.Range("K80").Value = 3424.9  # was 3694.3333
$ws.Range("L80").Value = 4870.5  # was 4836.727
$ws.Range("M80").Value = -2426.9  # was -2696.3333
$ws.Range("N80").Value = -6866.5  # was -6832.727

# Row 83
$ws.Range("H83").Value = 4147.7  # was 4322.65
$ws.Range("I83").Value = 3424.9  # was 3694.3333
$ws.Range("J83").Value = 4870.5  # was 4836.727
$ws.Range("K83").Value = 17124.5  # was 18471.6665
$ws.Range("L83").Value = 24352.5  # was 24183.635
$ws.Range("M83").Value = -12132.5  # was -13479.6665
$ws.Range("N83").Value = -34336.5  # was -34167.63499999999

# Row 122
$ws.Range("H122").Value = 85789.914  # was 65010.5
$ws.Range("I122").Value = 2160.375  # was 2097.8333
$ws.Range("J122").Value = 253049  # was 253748.5
$ws.Range("K122").Value = 6481.125  # was 6293.499899999999
$ws.Range("L122").Value = 759147  # was 761245.5
$ws.Range("M122").Value = -4031.125  # was -3843.499899999999
$ws.Range("N122").Value = -764047  # was -766145.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0  # was 700
$ws.Range("I7").Value = 0  # was 700
$ws.Range("K7").Value = 0  # was 700
$ws.Range("M7").ClearContents()  # was -588

# Row 40
$ws.Range("H40").Value = 4000  # was 0
$ws.Range("I40").Value = 4000  # was 0
$ws.Range("K40").Value = 4000  # was 0
$ws.Range("M40").Value = -3864  # was absent

# Row 46
$ws.Range("H46").Value = 80538.46000000001  # was 114333.445
$ws.Range("I46").Value = 3143  # was 3167
$ws.Range("J46").Value = 170833.17  # was 336666.34
$ws.Range("K46").Value = 3143  # was 3167
$ws.Range("L46").Value = 170833.17  # was 336666.34
$ws.Range("M46").Value = -2955  # was -2979
$ws.Range("N46").Value = -171209.17  # was -337042.34

# Row 55
$ws.Range("H55").Value = 1494.3334  # was 1495.75
$ws.Range("I55").Value = 993  # was 994
$ws.Range("J55").Value = 1995.6666  # was 1997.5
$ws.Range("K55").Value = 993  # was 994
$ws.Range("L55").Value = 1995.6666  # was 1997.5
$ws.Range("M55").Value = -820  # was -821
$ws.Range("N55").Value = -2341.6666  # was -2343.5

# Row 126
$ws.Range("H126").Value = 0  # was 700
$ws.Range("I126").Value = 0  # was 700
$ws.Range("K126").Value = 0  # was 2100
$ws.Range("M126").ClearContents()  # was 370

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 2489.158  # was 2461.55
$ws.Range("I100").Value = 2340.4614  # was 2311.6428
$ws.Range("K100").Value = 4680.9228  # was 4623.2856
$ws.Range("M100").Value = -4139.9228  # was -4082.2856

# Row 126
$ws.Range("H126").Value = 0  # was 3000
$ws.Range("I126").Value = 0  # was 3000
$ws.Range("K126").Value = 0  # was 9000
$ws.Range("M126").ClearContents()  # was -6530
